# DPUB-ARIA-test-results.xlsx edit
#
# Commit: "replaced 1 and 0 boolean values with the string exposed by
# developer inspector tool and accessibility insight tool"
#
# For every data row (4-41) on the "results" sheet, columns B
# ("builtin inspector ..."), C ("accessibility insights ..."), H
# ("inspector (Firefox ...)") and I ("accessibility insights (Firefox ...)")
# currently hold numeric booleans (1 = found / -1 = not found). They need to
# become the same text Excel already shows for a "found" result (the row's
# own doc-* role name, taken from column A) or the literal "-" placeholder
# that the sheet already uses elsewhere for "not found".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 4; $r -le 41; $r++) {
    $roleName = $ws.Cells.Item($r, 1).Value()

    foreach ($col in 2, 3, 8, 9) {
        $cell = $ws.Cells.Item($r, $col)
        $current = $cell.Value()

        if ($current -eq $null) {
            continue
        }

        if ($current -eq 1) {
            $cell.Value = $roleName
        } elseif ($current -eq -1) {
            $cell.Value = "-"
        }
    }
}

# Restore/point the selection at the cell the sheet was left on after the
# edits (bottom-right frozen pane).
$ws.Range("V41").Select()
